$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JewelleryPage")
$ws.Select()

# Clear any existing content within the range we will rewrite, just in case
$ws.Range("A1:C7").ClearContents()

# Write cells in a specific order so the shared-strings table is built
# in the same sequence as the target workbook.
$ws.Range("A1").Value = "MenuOption"
$ws.Range("B1").Value = "`n  SortOptions     "
$ws.Range("B2").Value = "Position"
$ws.Range("B3").Value = "Name: A to Z"
$ws.Range("B4").Value = "Name: Z to A"
$ws.Range("B5").Value = "Price: Low to High"
$ws.Range("C2").Value = "Black & White Diamond Heart"
$ws.Range("B6").Value = "Price: High to Low"
$ws.Range("B7").Value = "Created on"
$ws.Range("A2").Value = "JEWELRY"
$ws.Range("C1").Value = "productName"

# Formatting: wrap text on B1 and A2 (style index 2 in target)
$ws.Range("B1").WrapText = $true
$ws.Range("A2").WrapText = $true

# Row height for row 1
$ws.Rows.Item(1).RowHeight = 32.4

# Column widths
$ws.Columns.Item(1).ColumnWidth = 17.5546875
$ws.Columns.Item(2).ColumnWidth = 25.44140625
$ws.Columns.Item(3).ColumnWidth = 26.21875

# Selection state to match target
$ws.Range("B12").Select() | Out-Null
